# Project Starfighter 1.0 - apply edits:
#  1. Update two code-comment annotations on the "AnimatedSprite Test Code" slide.
#  2. Delete the two "Player Test Code" slides (Facing/Thrusting and ScrollRate/AccelRate).

$p = $ppt.ActivePresentation

# --- 1. Text edits on the "AnimatedSprite Test Code" slide (slide 18) ---
$animSlide = $p.Slides.Item(18)
$contentShape = $animSlide.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

$text = $tr.Text
$firstIdx = $text.IndexOf("[TestMethod]")
$bracketPos1 = $firstIdx + 12
$cr1 = $tr.Characters($bracketPos1, 1)
$cr1.Text = "] //sprite can update frames on its own"

$text = $tr.Text
$firstIdx = $text.IndexOf("[TestMethod]")
$secondIdx = $text.IndexOf("[TestMethod]", $firstIdx + 1)
$bracketPos2 = $secondIdx + 12
$cr2 = $tr.Characters($bracketPos2, 1)
$cr2.Text = "]// gets correct image from sprite sheet"

# --- 2. Delete the two "Player Test Code" slides ---
# (slide 27: Facing_Test/Thrusting_Test, slide 28: ScrollRate_Test/AccelRate_Test)
$p.Slides.Item(27).Delete()
$p.Slides.Item(27).Delete()
